# Auto-generated Excel COM-interop script to apply market-data refresh changes
# to the Ifrit_Profits workbook (per-sheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 18811.215
$ws.Range("I21").Value = 20105.7
$ws.Range("K21").Value = 20105.7
$ws.Range("M21").Value = -19637.7
$ws.Range("H23").Value = 18811.215
$ws.Range("I23").Value = 20105.7
$ws.Range("K23").Value = 20105.7
$ws.Range("M23").Value = -19871.7
$ws.Range("H29").Value = 800.75
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 78.63636
$ws.Range("I38").Value = 78.63636
$ws.Range("K38").Value = 235.90908
$ws.Range("M38").Value = 136.09092
$ws.Range("H64").Value = 3807.2727
$ws.Range("J64").Value = 3100
$ws.Range("L64").Value = 3100
$ws.Range("N64").Value = -3596
$ws.Range("H67").Value = 3807.2727
$ws.Range("J67").Value = 3100
$ws.Range("L67").Value = 3100
$ws.Range("N67").Value = -4816
$ws.Range("H98").Value = 3465.5
$ws.Range("I98").Value = 3103.8975
$ws.Range("J98").Value = 8166.3335
$ws.Range("K98").Value = 3103.8975
$ws.Range("L98").Value = 8166.3335
$ws.Range("M98").Value = -1605.8975
$ws.Range("N98").Value = -11162.3335
$ws.Range("H100").Value = 1474
$ws.Range("I100").Value = 1469.4445
$ws.Range("J100").Value = 1485.7142
$ws.Range("K100").Value = 1469.4445
$ws.Range("L100").Value = 1485.7142
$ws.Range("M100").Value = -928.4445000000001
$ws.Range("N100").Value = -2567.7142
$ws.Range("H112").Value = 41668068
$ws.Range("I112").Value = 699.75
$ws.Range("J112").Value = 50001540
$ws.Range("K112").Value = 2099.25
$ws.Range("L112").Value = 150004620
$ws.Range("M112").Value = -991.25
$ws.Range("N112").Value = -150006836
$ws.Range("H118").Value = 446.0909
$ws.Range("I118").Value = 276.55554
$ws.Range("K118").Value = 829.66662
$ws.Range("M118").Value = 827.33338
$ws.Range("H122").Value = 3465.5
$ws.Range("I122").Value = 3103.8975
$ws.Range("J122").Value = 8166.3335
$ws.Range("K122").Value = 9311.692500000001
$ws.Range("L122").Value = 24499.0005
$ws.Range("M122").Value = -6861.692500000001
$ws.Range("N122").Value = -29399.0005
$ws.Range("H129").Value = 966.5789
$ws.Range("I129").Value = 335.27274
$ws.Range("J129").Value = 1117.5435
$ws.Range("K129").Value = 1005.81822
$ws.Range("L129").Value = 3352.6305
$ws.Range("M129").Value = 3994.18178
$ws.Range("N129").Value = -13352.6305
$ws.Range("H138").Value = 3972.63
$ws.Range("I138").Value = 3632.3333
$ws.Range("J138").Value = 4019.0342
$ws.Range("K138").Value = 10896.9999
$ws.Range("L138").Value = 12057.1026
$ws.Range("M138").Value = -5756.999899999999
$ws.Range("N138").Value = -22337.1026
$ws.Range("H141").Value = 2409.5
$ws.Range("I141").Value = 1599.2858
$ws.Range("J141").Value = 4300
$ws.Range("K141").Value = 4797.857400000001
$ws.Range("L141").Value = 12900
$ws.Range("M141").Value = 382.1425999999992
$ws.Range("N141").Value = -23260

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 450.5
$ws.Range("I4").Value = 450.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 450.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -334.5
$ws.Range("N4").ClearContents()
$ws.Range("H44").Value = 22375
$ws.Range("J44").Value = 22375
$ws.Range("L44").Value = 22375
$ws.Range("N44").Value = -23351
$ws.Range("H117").Value = 20066.666
$ws.Range("J117").Value = 20066.666
$ws.Range("L117").Value = 20066.666
$ws.Range("N117").Value = -29244.666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6453394
$ws.Range("I134").Value = 7144654
$ws.Range("J134").Value = 1634.6666
$ws.Range("K134").Value = 21433962
$ws.Range("L134").Value = 4903.9998
$ws.Range("M134").Value = -21431427
$ws.Range("N134").Value = -9973.9998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23066.451
$ws.Range("I31").Value = 85910.336
$ws.Range("J31").Value = 7983.92
$ws.Range("K31").Value = 85910.336
$ws.Range("L31").Value = 7983.92
$ws.Range("M31").Value = -85615.336
$ws.Range("N31").Value = -8573.92
$ws.Range("H34").Value = 23066.451
$ws.Range("I34").Value = 85910.336
$ws.Range("J34").Value = 7983.92
$ws.Range("K34").Value = 85910.336
$ws.Range("L34").Value = 7983.92
$ws.Range("M34").Value = -85708.336
$ws.Range("N34").Value = -8387.92
$ws.Range("H50").Value = 11219
$ws.Range("J50").Value = 11219
$ws.Range("L50").Value = 11219
$ws.Range("N50").Value = -12469
$ws.Range("H51").Value = 10064
$ws.Range("J51").Value = 10425.143
$ws.Range("L51").Value = 10425.143
$ws.Range("N51").Value = -11897.143
$ws.Range("H60").Value = 9591.1
$ws.Range("I60").Value = 7750
$ws.Range("J60").Value = 10051.375
$ws.Range("K60").Value = 7750
$ws.Range("L60").Value = 10051.375
$ws.Range("M60").Value = -7239
$ws.Range("N60").Value = -11073.375
$ws.Range("H61").Value = 10064
$ws.Range("J61").Value = 10425.143
$ws.Range("L61").Value = 10425.143
$ws.Range("N61").Value = -11121.143
$ws.Range("H74").Value = 14421.1
$ws.Range("J74").Value = 16730.125
$ws.Range("L74").Value = 16730.125
$ws.Range("N74").Value = -18478.125
$ws.Range("H77").Value = 14421.1
$ws.Range("J77").Value = 16730.125
$ws.Range("L77").Value = 50190.375
$ws.Range("N77").Value = -58926.375
$ws.Range("H94").Value = 1547.0834
$ws.Range("I94").Value = 889.25
$ws.Range("K94").Value = 889.25
$ws.Range("M94").Value = -438.25
$ws.Range("H107").Value = 1895650.4
$ws.Range("I107").Value = 2316616.2
$ws.Range("J107").Value = 1303.25
$ws.Range("K107").Value = 2316616.2
$ws.Range("L107").Value = 1303.25
$ws.Range("M107").Value = -2314696.2
$ws.Range("N107").Value = -5143.25
$ws.Range("H132").Value = 1970.5454
$ws.Range("I132").Value = 1377.2
$ws.Range("J132").Value = 2465
$ws.Range("K132").Value = 4131.6
$ws.Range("L132").Value = 7395
$ws.Range("M132").Value = -1601.6
$ws.Range("N132").Value = -12455

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 370
$ws.Range("I92").Value = 375
$ws.Range("J92").Value = 366.66666
$ws.Range("K92").Value = 1125
$ws.Range("L92").Value = 1099.99998
$ws.Range("M92").Value = 123
$ws.Range("N92").Value = -3595.99998
$ws.Range("H100").Value = 2693.5833
$ws.Range("I100").Value = 2402.5
$ws.Range("J100").Value = 2751.8
$ws.Range("K100").Value = 7207.5
$ws.Range("L100").Value = 8255.400000000001
$ws.Range("M100").Value = -6396.5
$ws.Range("N100").Value = -9877.400000000001
$ws.Range("H106").Value = 4277537.5
$ws.Range("J106").Value = 4277537.5
$ws.Range("L106").Value = 12832612.5
$ws.Range("N106").Value = -12834504.5
$ws.Range("H114").Value = 1237.8823
$ws.Range("J114").Value = 1415.4348
$ws.Range("L114").Value = 4246.3044
$ws.Range("N114").Value = -10754.3044
$ws.Range("H131").Value = 48616320
$ws.Range("I131").Value = 133345480
$ws.Range("J131").Value = 26319174
$ws.Range("K131").Value = 400036440
$ws.Range("L131").Value = 78957522
$ws.Range("M131").Value = -400031400
$ws.Range("N131").Value = -78967602

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 102872.37
$ws.Range("I80").Value = 3615
$ws.Range("J80").Value = 159590.86
$ws.Range("K80").Value = 3615
$ws.Range("L80").Value = 159590.86
$ws.Range("M80").Value = -2617
$ws.Range("N80").Value = -161586.86
$ws.Range("H83").Value = 102872.37
$ws.Range("I83").Value = 3615
$ws.Range("J83").Value = 159590.86
$ws.Range("K83").Value = 18075
$ws.Range("L83").Value = 797954.2999999999
$ws.Range("M83").Value = -13083
$ws.Range("N83").Value = -807938.2999999999
$ws.Range("H122").Value = 4841.476
$ws.Range("I122").Value = 5851.615
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 17554.845
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -15104.845
$ws.Range("N122").Value = -14500
$ws.Range("H132").Value = 1816.4783
$ws.Range("I132").Value = 1285.4375
$ws.Range("J132").Value = 3030.2856
$ws.Range("K132").Value = 3856.3125
$ws.Range("L132").Value = 9090.856800000001
$ws.Range("M132").Value = -1326.3125
$ws.Range("N132").Value = -14150.8568

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1554.75
$ws.Range("I16").Value = 1618.4
$ws.Range("K16").Value = 1618.4
$ws.Range("M16").Value = -1448.4
$ws.Range("H46").Value = 1017.5455
$ws.Range("I46").Value = 1020.4737
$ws.Range("J46").Value = 999
$ws.Range("K46").Value = 1020.4737
$ws.Range("L46").Value = 999
$ws.Range("M46").Value = -832.4737
$ws.Range("N46").Value = -1375
$ws.Range("H55").Value = 122.71429
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 145.42857
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 145.42857
$ws.Range("M55").Value = 73
$ws.Range("N55").Value = -491.42857
$ws.Range("H136").Value = 2118.2942
$ws.Range("I136").Value = 1300.7142
$ws.Range("J136").Value = 5933.6665
$ws.Range("K136").Value = 3902.1426
$ws.Range("L136").Value = 17800.9995
$ws.Range("M136").Value = -1352.1426
$ws.Range("N136").Value = -22900.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3728.037
$ws.Range("I132").Value = 6450.5
$ws.Range("J132").Value = 2126.5881
$ws.Range("K132").Value = 19351.5
$ws.Range("L132").Value = 6379.7643
$ws.Range("M132").Value = -16821.5
$ws.Range("N132").Value = -11439.7643
$ws.Range("H136").Value = 14336.238
$ws.Range("I136").Value = 21502
$ws.Range("J136").Value = 2691.875
$ws.Range("K136").Value = 64506
$ws.Range("L136").Value = 8075.625
$ws.Range("M136").Value = -61956
$ws.Range("N136").Value = -13175.625
